# Helper: the COM host stores Shape.Left/Top/Width/Height internally as
# single-precision (float32) points and truncates back to EMU on save, so a
# naive "$emu / 12700" assignment can land 1 EMU short of the target. Nudge
# the point value up by the smallest amount needed so it round-trips to the
# exact target EMU.
function EmuToPoints($emu) {
    $basePts = $emu / 12700.0
    for ($i = 0; $i -lt 2000000; $i++) {
        $pts = $basePts + ($i * 0.0000001)
        $roundTripEmu = [math]::Floor([float]$pts * 12700)
        if ($roundTripEmu -eq $emu) {
            return $pts
        }
    }
    return $basePts
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# Shape 1 ("Rectangle 4"): move/resize, give it an opaque white fill, drop
# the line, and clear any effects.
$rect = $s.Shapes.Item(1)
$rect.Left   = EmuToPoints(1535456)
$rect.Top    = EmuToPoints(1310240)
$rect.Width  = EmuToPoints(5706775)
$rect.Height = EmuToPoints(2934061)

$rect.Fill.Visible = $true
$rect.Fill.Solid()
$rect.Fill.ForeColor.RGB = 16777215
$rect.Line.Visible = $false
$rect.Shadow.Visible = $false

# Shape 2 ("Picture 2" / Hitch - New Page.png): reposition only.
$pic = $s.Shapes.Item(2)
$pic.Left = EmuToPoints(1794475)
$pic.Top  = EmuToPoints(1226691)
